# ======================================================================
# Adds a new "2022-Q3" worksheet (positioned right after "总计") holding
# the quarterly fund-holdings breakdown, and records the 2022-Q3 summary
# row at the top of the "总计" sheet's data (existing rows shift down
# by one, which is also what surfaces the trailing 2020-Q4 total row).
# ======================================================================

$wb  = $excel.ActiveWorkbook
$tot = $wb.Worksheets.Item(1)          # "总计" — stays the first sheet

# ---- 1) new "2022-Q3" sheet, inserted right after "总计" --------------
$q3 = $wb.Worksheets.Add($null, $tot)
$q3.Name = "2022-Q3"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$hdrRange = $q3.Range("B1:H1")
$hdrArr = New-Object "object[,]" 1,7
for ($c = 0; $c -lt 7; $c++) { $hdrArr[0,$c] = $headers[$c] }
$hdrRange.Value = $hdrArr
$hdrRange.Font.Bold = $true
$hdrRange.Borders.LineStyle = 1
$hdrRange.HorizontalAlignment = -4108   # xlCenter
$hdrRange.VerticalAlignment = -4160     # xlTop

# Fund rows: A = 0-based index (number), B..G = text (fund code / name /
# scale / position / ratio / value all stored as text, matching the source
# data), H = rank (number).
$q3Data = @(
    @("011333","鹏华品质优选混合A","35.42","91.82","8.81","3.1205",4),
    @("100061","富国中国中小盘混合（QDII）人民币","35.11","83.32","3.28","1.1516",6),
    @("010591","富国中国中小盘混合（QDII）美元","35.11","83.32","3.28","1.1516",6),
    @("012057","鹏华品质成长混合A","11.11","84.05","8.88","0.9866",4),
    @("011570","鹏华鑫远价值一年持有期混合A","9.62","91.77","9.77","0.9399",2),
    @("008283","易方达金融行业股票","17.50","86.30","3.79","0.6632",10),
    @("009984","鹏华启航混合","12.45","84.87","4.80","0.5976",4),
    @("009234","鹏华优质企业混合","3.64","89.96","9.32","0.3392",3),
    @("005583","易方达港股通红利灵活配置混合","6.98","90.31","4.55","0.3176",4),
    @("009223","宝盈现代服务业混合A","3.18","93.59","7.69","0.2445",2),
    @("011334","鹏华品质优选混合C","2.41","91.82","8.81","0.2123",4),
    @("013859","宝盈品质甄选混合A","1.80","93.72","7.74","0.1393",2),
    @("013334","鹏华价值远航6个月持有期混合A","1.56","90.27","8.06","0.1257",4),
    @("001703","银华沪港深增长股票A","2.07","88.44","5.07","0.1049",4),
    @("501021","华宝标普香港上市中国中小盘指数（LOF）A","4.19","92.99","2.39","0.1001",2),
    @("012640","鹏华稳健鸿利一年持有期混合A","2.61","92.98","3.43","0.0895",7),
    @("008134","鹏华优选价值股票","1.80","92.72","3.37","0.0607",8),
    @("011571","鹏华鑫远价值一年持有期混合C","0.59","91.77","9.77","0.0576",2),
    @("008303","宝盈龙头优选股票A","0.65","93.98","7.71","0.0501",2),
    @("012058","鹏华品质成长混合C","0.39","84.05","8.88","0.0346",4),
    @("009224","宝盈现代服务业混合C","0.43","93.59","7.69","0.0331",2),
    @("011969","建信港股通精选混合A","0.56","63.37","3.81","0.0213",10),
    @("006675","宝盈品牌消费股票A","0.17","93.83","7.71","0.0131",2),
    @("014364","银华沪港深增长股票C","0.25","88.44","5.07","0.0127",4),
    @("008304","宝盈龙头优选股票C","0.16","93.98","7.71","0.0123",2),
    @("006676","宝盈品牌消费股票C","0.14","93.83","7.71","0.0108",2),
    @("011970","建信港股通精选混合C","0.23","63.37","3.81","0.0088",10),
    @("013335","鹏华价值远航6个月持有期混合C","0.10","90.27","8.06","0.0081",4),
    @("378006","上投摩根全球新兴市场混合（QDII）","0.40","87.48","1.91","0.0076",9),
    @("004532","民生加银中证港股通高股息精选指数A","0.13","92.87","5.80","0.0075",1),
    @("011647","博时港股通红利精选混合A","0.11","82.44","6.25","0.0069",4),
    @("006127","华宝标普香港上市中国中小盘指数（LOF）C","0.24","92.99","2.39","0.0057",2),
    @("013860","宝盈品质甄选混合C","0.07","93.72","7.74","0.0054",2),
    @("501303","广发恒生中型股指数（LOF）A","0.21","89.12","2.43","0.0051",2),
    @("006658","财通中证香港红利等权投资指数A","0.13","88.79","3.86","0.0050",5),
    @("004533","民生加银中证港股通高股息精选指数C","0.08","92.87","5.80","0.0046",1),
    @("012641","鹏华稳健鸿利一年持有期混合C","0.10","92.98","3.43","0.0034",7),
    @("004996","广发恒生中型股指数（LOF）C","0.09","89.12","2.43","0.0022",2),
    @("160922","大成恒生综合中小型股指数（QDII-LOF）A","0.09","86.62","1.74","0.0016",2),
    @("006659","财通中证香港红利等权投资指数C","0.04","88.79","3.86","0.0015",5),
    @("011648","博时港股通红利精选混合C","0.02","82.44","6.25","0.0012",4),
    @("008972","大成恒生综合中小型股指数C","0.02","86.62","1.74","0.0003",2)
)

$n = $q3Data.Count

$idxRange = $q3.Range("A2:A" + (1 + $n))
$idxArr = New-Object "object[,]" $n,1
for ($r = 0; $r -lt $n; $r++) { $idxArr[$r,0] = $r }
$idxRange.Value = $idxArr
$idxRange.Font.Bold = $true
$idxRange.Borders.LineStyle = 1
$idxRange.HorizontalAlignment = -4108
$idxRange.VerticalAlignment = -4160

$bodyRange = $q3.Range("B2:G" + (1 + $n))
$bodyRange.NumberFormat = "@"        # keep fund codes/percentages as text
$bodyArr = New-Object "object[,]" $n,6
for ($r = 0; $r -lt $n; $r++) {
    for ($c = 0; $c -lt 6; $c++) {
        $bodyArr[$r,$c] = $q3Data[$r][$c]
    }
}
$bodyRange.Value = $bodyArr

$rankRange = $q3.Range("H2:H" + (1 + $n))
$rankArr = New-Object "object[,]" $n,1
for ($r = 0; $r -lt $n; $r++) { $rankArr[$r,0] = $q3Data[$r][6] }
$rankRange.Value = $rankArr

$q3.Range("A1").Select()

# ---- 2) "总计": insert the 2022-Q3 summary row at the top of the data -
$tot.Rows.Item(2).Insert()

$tot.Range("A2").Value = 0
$tot.Range("A2").Font.Bold = $true
$tot.Range("A2").Borders.LineStyle = 1
$tot.Range("A2").HorizontalAlignment = -4108
$tot.Range("A2").VerticalAlignment = -4160
$tot.Range("B2").Value = "2022-Q3"
$tot.Range("C2").Value = 42
$tot.Range("D2").Value = 10.67

$tot.Range("A1").Select()
